$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-6 from 45175 to 45183
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45183
}
